# Updated Masterdata as per 2nd may Data Refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user")

# Refresh the regcntr_id values that changed in the latest DB extract.
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Re-select/scroll the sheet to where the editor left off when saving.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C19").Select()
